# added a warning to popup
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining data row: the "Day" moved to 09.08.2022 and the
# desk assignment to Desk B35. Force the Day cell to stay text (it looks
# like a date and would otherwise get auto-converted to a date serial).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "09.08.2022"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "Desk B35"

# Remove the old duplicate entries in rows 3-8, leaving only the header
# and the single remaining row.
$ws.Range("A3:D8").EntireRow.Delete()
